$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer readings to insert right after the header row (row 1),
# pushing all existing data rows down by 11 rows.
$newData = @(
    @(-1.931140422821045, 9.274446487426758, 0.347445011138916),
    @(-2.475548833608627, 9.384642362594603, 0.7847917079925537),
    @(-2.621871948242188, 9.314098954200743, 1.41282993555069),
    @(-2.868601083755493, 9.443870902061462, 0.9555243626236917),
    @(-2.620113015174865, 9.546792268753052, 0.6990440487861633),
    @(-2.44504114985466, 9.541788041591644, 0.372002582065761),
    @(-2.28140389919281, 9.524857640266418, -0.01871592737734312),
    @(-2.462455779314041, 9.538427114486694, -0.03734804317355161),
    @(-2.66546654701233, 9.521270275115967, 0.1833332777023315),
    @(-2.735388696193695, 9.501047194004059, 0.2642159881070256),
    @(-2.759680032730102, 9.431608200073242, 0.08809284307062609)
)

$rowCount = $newData.Count

# Insert the required number of blank rows right after the header row (row 2..12),
# pushing all the pre-existing data rows down by $rowCount rows.
$insertRange = $ws.Range("A2:A$($rowCount + 1)").EntireRow
$insertRange.Insert()

# The inserted rows pick up the header row's formatting by default; the
# original data rows carry no explicit style, so clear it back off.
$ws.Range("A2:C$($rowCount + 1)").ClearFormats()

# Fill in the new rows with data
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
    $ws.Cells.Item($r, 3).Value = $newData[$i][2]
}

# The old last data row (originally row 21) is dropped from the data set.
# After the insertion above it now lives at row 21 + $rowCount.
$oldLastRow = 21 + $rowCount
$ws.Range("A$($oldLastRow):C$($oldLastRow)").EntireRow.Delete()
